$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Data for the three 10-day blocks (rows 200-209, 218-227, 236-245).
# Each block: B/C/D raw 0|1 values, and an optional "megjegyzes" (E column)
# shared-string comment.
# ---------------------------------------------------------------------------

$block1 = @(
  @{ Row=200; B=1; C=1; D=1; E=$null },
  @{ Row=201; B=1; C=1; D=1; E="4sz" },
  @{ Row=202; B=0; C=1; D=1; E=$null },
  @{ Row=203; B=1; C=0; D=1; E=$null },
  @{ Row=204; B=1; C=1; D=1; E="2sz" },
  @{ Row=205; B=1; C=1; D=0; E="2sz" },
  @{ Row=206; B=1; C=0; D=0; E=$null },
  @{ Row=207; B=1; C=1; D=0; E=$null },
  @{ Row=208; B=0; C=0; D=1; E=$null },
  @{ Row=209; B=0; C=0; D=0; E=$null }
)

$block2 = @(
  @{ Row=218; B=1; C=1; D=1; E="2e 2sz" },
  @{ Row=219; B=0; C=0; D=0; E=$null },
  @{ Row=220; B=1; C=0; D=1; E=$null },
  @{ Row=221; B=0; C=1; D=1; E=$null },
  @{ Row=222; B=1; C=1; D=1; E="2nsz" },
  @{ Row=223; B=1; C=1; D=0; E="3sz" },
  @{ Row=224; B=1; C=0; D=0; E=$null },
  @{ Row=225; B=1; C=0; D=0; E=$null },
  @{ Row=226; B=1; C=1; D=0; E=$null },
  @{ Row=227; B=1; C=0; D=1; E="3e" }
)

$block3 = @(
  @{ Row=236; B=1; C=1; D=0; E="2sz" },
  @{ Row=237; B=1; C=1; D=1; E="4sz" },
  @{ Row=238; B=1; C=1; D=1; E=$null },
  @{ Row=239; B=1; C=1; D=1; E=$null },
  @{ Row=240; B=1; C=1; D=1; E=$null },
  @{ Row=241; B=1; C=1; D=1; E="2sz" },
  @{ Row=242; B=0; C=0; D=0; E=$null },
  @{ Row=243; B=1; C=0; D=1; E=$null },
  @{ Row=244; B=1; C=1; D=1; E="2sz" },
  @{ Row=245; B=1; C=0; D=0; E=$null }
)

$blocks = @($block1, $block2, $block3)

foreach ($blk in $blocks) {
  foreach ($item in $blk) {
    $r = $item.Row
    $ws.Range("B$r").Value = $item.B
    $ws.Range("C$r").Value = $item.C
    $ws.Range("D$r").Value = $item.D
    if ($item.E -ne $null) {
      $ws.Range("E$r").Value = $item.E
    }
  }

  # First row of the block gets its own (non-shared) AVERAGE formula.
  $firstRow = $blk[0].Row
  $ws.Range("G$firstRow").Formula = "=AVERAGE(B$firstRow`:D$firstRow)"

  # The remaining nine rows share one formula definition, entered on the
  # whole range at once so the engine emits a shared formula (t="shared").
  $restStart = $blk[1].Row
  $restEnd = $blk[$blk.Count - 1].Row
  $ws.Range("G$restStart`:G$restEnd").Formula = "=AVERAGE(B$restStart`:D$restStart)"
}

# ---------------------------------------------------------------------------
# Day-average formulas in O22:O24 (one per 18-day block).
# ---------------------------------------------------------------------------
$ws.Range("O22").Formula = "=AVERAGE(B200:B217)"
$ws.Range("O23").Formula = "=AVERAGE(B218:B235)"
$ws.Range("O24").Formula = "=AVERAGE(B236:B253)"

# ---------------------------------------------------------------------------
# Sheet view: scroll back to the top (drop topLeftCell) and move the active
# selection to O25.
# ---------------------------------------------------------------------------
[void]$ws.Range("A1").Select()
[void]$ws.Range("O25").Select()
